# Weekly driver report update for 2025-04-20
# Updates the "Bad Drivers" summary stats and refreshes the "Good Drivers"
# table with this week's driver vintage / sample-count data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers summary (rows 3-5) ---------------------------------------
$ws.Range("D3").Value = 97.7

$ws.Range("C4").Value = 164
$ws.Range("D4").Value = 98.09999999999999

$ws.Range("C5").Value = 243

# --- Good Drivers table (rows 13-18) ---------------------------------------
# Row 13: Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B13").Value = 445055
$ws.Range("D13").Value = 99.90000000000001
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2024-11-10"

# Row 14: Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B14").Value = 77849
$ws.Range("D14").Value = 99.90000000000001
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2021-08-18"

# Row 15: Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B15").Value = 34244
$ws.Range("D15").Value = 100
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2021-04-27"

# Row 16: Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B16").Value = 59673
$ws.Range("D16").Value = 100
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2020-08-05"

# Row 17: Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B17").Value = 113652
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2020-01-06"

# Row 18: Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1
$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B18").Value = 56018
